# TRIUNFO.xlsx update
# - Rename "Paineis DARQ" -> "PAINEIS DARQ"
# - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

# Suppress the "permanently delete" confirmation dialog that Excel would
# otherwise show when a sheet is removed via COM automation.
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
